{"js": "// The document contains a 5-column table of two-digit multiplication\n// facts (\"A\u00d7B=C\"), one equation per table cell run. The edit swaps each\n// equation for a new one. Every \"old\" equation string occurs exactly\n// once in the document body, and no \"old\" string is a substring of\n// another \"old\"/\"new\" string, so a plain search-and-replace per pair is\n// safe, unambiguous, and independent of iteration order.\nconst pairs = [\n  [\"68\u00d716=1088\", \"39\u00d737=1443\"],\n  [\"77\u00d710=770\", \"57\u00d773=4161\"],\n  [\"79\u00d713=1027\", \"65\u00d777=5005\"],\n  [\"84\u00d773=6132\", \"49\u00d723=1127\"],\n  [\"67\u00d712=804\", \"91\u00d754=4914\"],\n  [\"12\u00d735=420\", \"55\u00d794=5170\"],\n  [\"93\u00d798=9114\", \"20\u00d789=1780\"],\n  [\"79\u00d769=5451\", \"89\u00d739=3471\"],\n  [\"19\u00d762=1178\", \"95\u00d791=8645\"],\n  [\"21\u00d768=1428\", \"99\u00d771=7029\"],\n  [\"33\u00d767=2211\", \"55\u00d759=3245\"],\n  [\"85\u00d719=1615\", \"62\u00d759=3658\"],\n  [\"49\u00d768=3332\", \"85\u00d7100=8500\"],\n  [\"81\u00d728=2268\", \"14\u00d771=994\"],\n  [\"74\u00d763=4662\", \"31\u00d797=3007\"],\n  [\"59\u00d749=2891\", \"23\u00d783=1909\"],\n  [\"32\u00d730=960\", \"63\u00d755=3465\"],\n  [\"35\u00d767=2345\", \"51\u00d794=4794\"],\n  [\"11\u00d792=1012\", \"58\u00d739=2262\"],\n  [\"95\u00d774=7030\", \"87\u00d743=3741\"],\n  [\"25\u00d761=1525\", \"12\u00d714=168\"],\n  [\"12\u00d716=192\", \"29\u00d782=2378\"],\n  [\"16\u00d761=976\", \"18\u00d734=612\"],\n  [\"57\u00d744=2508\", \"25\u00d744=1100\"],\n  [\"36\u00d741=1476\", \"34\u00d760=2040\"],\n  [\"90\u00d777=6930\", \"98\u00d794=9212\"],\n  [\"52\u00d748=2496\", \"67\u00d764=4288\"],\n  [\"10\u00d757=570\", \"59\u00d777=4543\"],\n  [\"80\u00d710=800\", \"13\u00d773=949\"],\n  [\"15\u00d790=1350\", \"27\u00d746=1242\"],\n  [\"63\u00d799=6237\", \"77\u00d746=3542\"],\n  [\"51\u00d751=2601\", \"24\u00d713=312\"],\n  [\"46\u00d792=4232\", \"16\u00d777=1232\"],\n  [\"12\u00d721=252\", \"68\u00d773=4964\"],\n  [\"85\u00d785=7225\", \"83\u00d740=3320\"],\n  [\"14\u00d753=742\", \"34\u00d774=2516\"],\n  [\"21\u00d764=1344\", \"62\u00d760=3720\"],\n  [\"57\u00d751=2907\", \"79\u00d758=4582\"],\n  [\"28\u00d744=1232\", \"47\u00d737=1739\"],\n  [\"10\u00d790=900\", \"56\u00d745=2520\"],\n  [\"25\u00d784=2100\", \"18\u00d723=414\"],\n  [\"56\u00d722=1232\", \"71\u00d742=2982\"],\n  [\"35\u00d724=840\", \"62\u00d753=3286\"],\n  [\"100\u00d746=4600\", \"31\u00d727=837\"],\n  [\"21\u00d756=1176\", \"53\u00d764=3392\"],\n  [\"63\u00d797=6111\", \"10\u00d789=890\"],\n  [\"13\u00d781=1053\", \"88\u00d769=6072\"],\n  [\"43\u00d753=2279\", \"69\u00d788=6072\"],\n  [\"32\u00d787=2784\", \"68\u00d760=4080\"],\n  [\"46\u00d736=1656\", \"35\u00d769=2415\"],\n  [\"54\u00d750=2700\", \"30\u00d748=1440\"],\n  [\"38\u00d776=2888\", \"82\u00d728=2296\"],\n  [\"100\u00d792=9200\", \"24\u00d740=960\"],\n  [\"30\u00d771=2130\", \"13\u00d790=1170\"],\n  [\"45\u00d755=2475\", \"52\u00d788=4576\"],\n  [\"49\u00d736=1764\", \"79\u00d725=1975\"],\n  [\"76\u00d791=6916\", \"63\u00d736=2268\"],\n  [\"76\u00d766=5016\", \"89\u00d769=6141\"],\n  [\"11\u00d788=968\", \"42\u00d728=1176\"],\n  [\"20\u00d773=1460\", \"63\u00d781=5103\"],\n  [\"54\u00d719=1026\", \"54\u00d775=4050\"],\n  [\"82\u00d760=4920\", \"88\u00d715=1320\"],\n  [\"72\u00d728=2016\", \"71\u00d761=4331\"],\n  [\"26\u00d724=624\", \"90\u00d767=6030\"],\n  [\"58\u00d749=2842\", \"44\u00d791=4004\"],\n  [\"90\u00d712=1080\", \"94\u00d744=4136\"],\n  [\"51\u00d739=1989\", \"19\u00d781=1539\"],\n  [\"77\u00d774=5698\", \"96\u00d713=1248\"],\n  [\"28\u00d779=2212\", \"62\u00d721=1302\"],\n  [\"69\u00d727=1863\", \"32\u00d784=2688\"],\n  [\"56\u00d724=1344\", \"47\u00d768=3196\"],\n  [\"11\u00d732=352\", \"100\u00d797=9700\"],\n  [\"11\u00d750=550\", \"64\u00d721=1344\"],\n  [\"95\u00d772=6840\", \"23\u00d739=897\"],\n  [\"72\u00d787=6264\", \"85\u00d764=5440\"],\n  [\"28\u00d780=2240\", \"16\u00d729=464\"],\n  [\"58\u00d747=2726\", \"17\u00d750=850\"],\n  [\"93\u00d776=7068\", \"44\u00d778=3432\"],\n  [\"76\u00d799=7524\", \"90\u00d740=3600\"],\n  [\"29\u00d746=1334\", \"33\u00d788=2904\"],\n  [\"10\u00d773=730\", \"96\u00d786=8256\"],\n  [\"62\u00d776=4712\", \"54\u00d718=972\"],\n  [\"50\u00d721=1050\", \"36\u00d796=3456\"],\n  [\"71\u00d715=1065\", \"58\u00d721=1218\"],\n  [\"75\u00d778=5850\", \"92\u00d754=4968\"],\n  [\"90\u00d747=4230\", \"42\u00d763=2646\"],\n  [\"84\u00d728=2352\", \"79\u00d710=790\"],\n  [\"86\u00d738=3268\", \"76\u00d731=2356\"],\n  [\"50\u00d720=1000\", \"16\u00d785=1360\"],\n  [\"63\u00d796=6048\", \"74\u00d787=6438\"],\n  [\"65\u00d734=2210\", \"61\u00d792=5612\"],\n  [\"32\u00d759=1888\", \"27\u00d757=1539\"],\n  [\"92\u00d721=1932\", \"25\u00d777=1925\"],\n  [\"32\u00d785=2720\", \"72\u00d769=4968\"],\n  [\"89\u00d772=6408\", \"100\u00d768=6800\"],\n  [\"73\u00d775=5475\", \"66\u00d742=2772\"],\n  [\"34\u00d767=2278\", \"22\u00d764=1408\"],\n  [\"42\u00d794=3948\", \"67\u00d725=1675\"],\n  [\"26\u00d762=1612\", \"84\u00d736=3024\"],\n  [\"67\u00d775=5025\", \"97\u00d738=3686\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each two-digit multiplication equation with its new value.\n# Every \"old\" string occurs exactly once in the document body (inside a\n# table cell), so a plain Find/Replace-all per pair is unambiguous and\n# order-independent.\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"68\u00d716=1088\", \"39\u00d737=1443\"),\n    @(\"77\u00d710=770\", \"57\u00d773=4161\"),\n    @(\"79\u00d713=1027\", \"65\u00d777=5005\"),\n    @(\"84\u00d773=6132\", \"49\u00d723=1127\"),\n    @(\"67\u00d712=804\", \"91\u00d754=4914\"),\n    @(\"12\u00d735=420\", \"55\u00d794=5170\"),\n    @(\"93\u00d798=9114\", \"20\u00d789=1780\"),\n    @(\"79\u00d769=5451\", \"89\u00d739=3471\"),\n    @(\"19\u00d762=1178\", \"95\u00d791=8645\"),\n    @(\"21\u00d768=1428\", \"99\u00d771=7029\"),\n    @(\"33\u00d767=2211\", \"55\u00d759=3245\"),\n    @(\"85\u00d719=1615\", \"62\u00d759=3658\"),\n    @(\"49\u00d768=3332\", \"85\u00d7100=8500\"),\n    @(\"81\u00d728=2268\", \"14\u00d771=994\"),\n    @(\"74\u00d763=4662\", \"31\u00d797=3007\"),\n    @(\"59\u00d749=2891\", \"23\u00d783=1909\"),\n    @(\"32\u00d730=960\", \"63\u00d755=3465\"),\n    @(\"35\u00d767=2345\", \"51\u00d794=4794\"),\n    @(\"11\u00d792=1012\", \"58\u00d739=2262\"),\n    @(\"95\u00d774=7030\", \"87\u00d743=3741\"),\n    @(\"25\u00d761=1525\", \"12\u00d714=168\"),\n    @(\"12\u00d716=192\", \"29\u00d782=2378\"),\n    @(\"16\u00d761=976\", \"18\u00d734=612\"),\n    @(\"57\u00d744=2508\", \"25\u00d744=1100\"),\n    @(\"36\u00d741=1476\", \"34\u00d760=2040\"),\n    @(\"90\u00d777=6930\", \"98\u00d794=9212\"),\n    @(\"52\u00d748=2496\", \"67\u00d764=4288\"),\n    @(\"10\u00d757=570\", \"59\u00d777=4543\"),\n    @(\"80\u00d710=800\", \"13\u00d773=949\"),\n    @(\"15\u00d790=1350\", \"27\u00d746=1242\"),\n    @(\"63\u00d799=6237\", \"77\u00d746=3542\"),\n    @(\"51\u00d751=2601\", \"24\u00d713=312\"),\n    @(\"46\u00d792=4232\", \"16\u00d777=1232\"),\n    @(\"12\u00d721=252\", \"68\u00d773=4964\"),\n    @(\"85\u00d785=7225\", \"83\u00d740=3320\"),\n    @(\"14\u00d753=742\", \"34\u00d774=2516\"),\n    @(\"21\u00d764=1344\", \"62\u00d760=3720\"),\n    @(\"57\u00d751=2907\", \"79\u00d758=4582\"),\n    @(\"28\u00d744=1232\", \"47\u00d737=1739\"),\n    @(\"10\u00d790=900\", \"56\u00d745=2520\"),\n    @(\"25\u00d784=2100\", \"18\u00d723=414\"),\n    @(\"56\u00d722=1232\", \"71\u00d742=2982\"),\n    @(\"35\u00d724=840\", \"62\u00d753=3286\"),\n    @(\"100\u00d746=4600\", \"31\u00d727=837\"),\n    @(\"21\u00d756=1176\", \"53\u00d764=3392\"),\n    @(\"63\u00d797=6111\", \"10\u00d789=890\"),\n    @(\"13\u00d781=1053\", \"88\u00d769=6072\"),\n    @(\"43\u00d753=2279\", \"69\u00d788=6072\"),\n    @(\"32\u00d787=2784\", \"68\u00d760=4080\"),\n    @(\"46\u00d736=1656\", \"35\u00d769=2415\"),\n    @(\"54\u00d750=2700\", \"30\u00d748=1440\"),\n    @(\"38\u00d776=2888\", \"82\u00d728=2296\"),\n    @(\"100\u00d792=9200\", \"24\u00d740=960\"),\n    @(\"30\u00d771=2130\", \"13\u00d790=1170\"),\n    @(\"45\u00d755=2475\", \"52\u00d788=4576\"),\n    @(\"49\u00d736=1764\", \"79\u00d725=1975\"),\n    @(\"76\u00d791=6916\", \"63\u00d736=2268\"),\n    @(\"76\u00d766=5016\", \"89\u00d769=6141\"),\n    @(\"11\u00d788=968\", \"42\u00d728=1176\"),\n    @(\"20\u00d773=1460\", \"63\u00d781=5103\"),\n    @(\"54\u00d719=1026\", \"54\u00d775=4050\"),\n    @(\"82\u00d760=4920\", \"88\u00d715=1320\"),\n    @(\"72\u00d728=2016\", \"71\u00d761=4331\"),\n    @(\"26\u00d724=624\", \"90\u00d767=6030\"),\n    @(\"58\u00d749=2842\", \"44\u00d791=4004\"),\n    @(\"90\u00d712=1080\", \"94\u00d744=4136\"),\n    @(\"51\u00d739=1989\", \"19\u00d781=1539\"),\n    @(\"77\u00d774=5698\", \"96\u00d713=1248\"),\n    @(\"28\u00d779=2212\", \"62\u00d721=1302\"),\n    @(\"69\u00d727=1863\", \"32\u00d784=2688\"),\n    @(\"56\u00d724=1344\", \"47\u00d768=3196\"),\n    @(\"11\u00d732=352\", \"100\u00d797=9700\"),\n    @(\"11\u00d750=550\", \"64\u00d721=1344\"),\n    @(\"95\u00d772=6840\", \"23\u00d739=897\"),\n    @(\"72\u00d787=6264\", \"85\u00d764=5440\"),\n    @(\"28\u00d780=2240\", \"16\u00d729=464\"),\n    @(\"58\u00d747=2726\", \"17\u00d750=850\"),\n    @(\"93\u00d776=7068\", \"44\u00d778=3432\"),\n    @(\"76\u00d799=7524\", \"90\u00d740=3600\"),\n    @(\"29\u00d746=1334\", \"33\u00d788=2904\"),\n    @(\"10\u00d773=730\", \"96\u00d786=8256\"),\n    @(\"62\u00d776=4712\", \"54\u00d718=972\"),\n    @(\"50\u00d721=1050\", \"36\u00d796=3456\"),\n    @(\"71\u00d715=1065\", \"58\u00d721=1218\"),\n    @(\"75\u00d778=5850\", \"92\u00d754=4968\"),\n    @(\"90\u00d747=4230\", \"42\u00d763=2646\"),\n    @(\"84\u00d728=2352\", \"79\u00d710=790\"),\n    @(\"86\u00d738=3268\", \"76\u00d731=2356\"),\n    @(\"50\u00d720=1000\", \"16\u00d785=1360\"),\n    @(\"63\u00d796=6048\", \"74\u00d787=6438\"),\n    @(\"65\u00d734=2210\", \"61\u00d792=5612\"),\n    @(\"32\u00d759=1888\", \"27\u00d757=1539\"),\n    @(\"92\u00d721=1932\", \"25\u00d777=1925\"),\n    @(\"32\u00d785=2720\", \"72\u00d769=4968\"),\n    @(\"89\u00d772=6408\", \"100\u00d768=6800\"),\n    @(\"73\u00d775=5475\", \"66\u00d742=2772\"),\n    @(\"34\u00d767=2278\", \"22\u00d764=1408\"),\n    @(\"42\u00d794=3948\", \"67\u00d725=1675\"),\n    @(\"26\u00d762=1612\", \"84\u00d736=3024\"),\n    @(\"67\u00d775=5025\", \"97\u00d738=3686\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll) | Out-Null\n}\n"}
